# Add a new "footer" paragraph right after the existing "Header" paragraph,
# at the very end of the document body.

$d = $word.ActiveDocument

# Collapse to the end of the document and insert a brand-new paragraph mark there.
$endRange = $d.Content
$endRange.Collapse(0)              # wdCollapseEnd = 0
$endRange.InsertParagraphAfter()

# The newly created paragraph is now the last paragraph in the document;
# fill it with the footer text.
$newPara = $d.Paragraphs.Last
$newPara.Range.InsertBefore("footer")
